$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.763.26"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'1.000"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "'303.45"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "'0.3812"
$ws.Range("E7").Value = "  +0.99%  "
$ws.Range("D8").Value = "'0.3631"
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "'51.20"
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("D10").Value = "'1.257"
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("D11").Value = "'0.08227"
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "'22.73"
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("D14").Value = "'6.548"
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").Value = "'7.486"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("D16").Value = "'0.00001241"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "1.655.94"
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").Value = "'97.87"
$ws.Range("E18").Value = "  +2.96%  "
$ws.Range("D19").Value = "'0.06985"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").Value = "'6.809"
$ws.Range("E20").Value = "  +3.58%  "
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").Value = "'12.82"
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("D24").Value = "23.762.73"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("D25").Value = "'2.539"
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("D26").Value = "'3.079"
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("D28").Value = "'151.17"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").Value = "'5.249"
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("D30").Value = "'134.51"
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("D31").Value = "1.841.34"
$ws.Range("E31").Value = "  +1.79%  "
$ws.Range("D32").Value = "'6.966"
$ws.Range("E32").Value = "  +4.60%  "
$ws.Range("D33").Value = "'2.277"
$ws.Range("E33").Value = "  +5.34%  "
$ws.Range("D34").Value = "'1.074"
$ws.Range("E34").Value = "  +1.92%  "
$ws.Range("D35").Value = "'11.93"
$ws.Range("E35").Value = "  +4.32%  "
$ws.Range("D36").Value = "'0.02831"
$ws.Range("E36").Value = "  +2.04%  "
$ws.Range("D37").Value = "'0.2526"
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("D38").Value = "'6.142"
$ws.Range("E38").Value = "  +2.02%  "
$ws.Range("D39").Value = "'0.08843"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("D40").Value = "'0.07103"
$ws.Range("E40").Value = "  -0.73%  "
$ws.Range("D41").Value = "'13.26"
$ws.Range("E41").Value = "  +9.05%  "
$ws.Range("D42").Value = "'0.7074"
$ws.Range("E42").Value = "  +0.96%  "
$ws.Range("D43").Value = "'1.349"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").Value = "'16.04"
$ws.Range("E44").Value = "  +1.10%  "
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("D46").Value = "'2.340"
$ws.Range("E46").Value = "  +2.47%  "
$ws.Range("D47").Value = "'1.000"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").Value = "'3.959"
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("D49").Value = "'0.07962"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("D50").Value = "'127.96"
$ws.Range("E50").Value = "  +0.79%  "
$ws.Range("D51").Value = "'1.195"
$ws.Range("E51").Value = "  -0.13%  "
